# Phase 9: Regenerated Served Documents (HTML) & Updated Walkthrough
#
# Replaces every "GSM-R" occurrence with "TETRA" (covers the
# "TETRA-GSM-R" -> "TETRA-TETRA", "Comunicaciones GSM-R" -> "Comunicaciones TETRA",
# "Protocolo GSM-R" -> "Protocolo TETRA" and "TETRA + GSM-R" -> "TETRA + TETRA"
# cases), and revises the availability figure from 99.95% to 99.5%
# everywhere it is mentioned in the document body (narrative text,
# table cells and list items).

$d = $word.ActiveDocument

# 1) GSM-R -> TETRA (communications protocol naming correction)
$d.Content.Find.Execute(
    "GSM-R", $true, $false, $false, $false, $false,
    $true, 1, $false, "TETRA", 2
) | Out-Null

# 2) 99.95% -> 99.5% (revised availability target, all occurrences:
#    summary bullet, heading "¿Por qué disponibilidad 99.95%?", table
#    cells, acceptance-criteria bullet, contractual-compliance text and
#    KPI table)
$d.Content.Find.Execute(
    "99.95%", $true, $false, $false, $false, $false,
    $true, 1, $false, "99.5%", 2
) | Out-Null
